# Add an "adductName" column (with value "[M-H]-" for every data row) to the
# "Corrected" sheet, as a new column C (shifting the existing blank/tissue
# columns one slot to the right), and make "Corrected" the active/selected
# sheet instead of "Original".

$wb = $excel.ActiveWorkbook

$original = $wb.Worksheets.Item("Original")
$corrected = $wb.Worksheets.Item("Corrected")

# Insert a new blank column before column C ("blank" header) on the
# "Corrected" sheet; everything from C onward shifts right to D onward.
$corrected.Columns("C:C").Insert()

# New header cell for the inserted column.
$corrected.Range("C1").Style = "Normal"
$corrected.Range("C1").Value = "adductName"
$corrected.Range("C1").Font.Bold = $true

# Fill the adduct name for every data row (rows 2-12).
for ($r = 2; $r -le 12; $r++) {
    $cell = $corrected.Cells.Item($r, 3)
    $cell.Style = "Normal"
    $cell.Value = "[M-H]-"
}

# "Corrected" becomes the active/selected sheet (was "Original").
$corrected.Activate()
$corrected.Select()
